$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nbsp = [char]0x00A0

$data = @(
    ,@(" Dubai (DSC)"," October 04 2020","Super Kings won by 10 wickets (with 14 balls remaining)","Chennai Super Kings","Kings XI Punjab","87","53","11","1","164.15")
    ,@(" Dubai (DSC)"," October 13 2020","Super Kings won by 20 runs","Chennai Super Kings","Sunrisers Hyderabad","0","1","0","0","0.00")
    ,@(" Abu Dhabi"," September 19 2020","Super Kings won by 5 wickets (with 4 balls remaining)","Chennai Super Kings","Mumbai Indians","58","44","6","0","131.81")
    ,@(" Dubai (DSC)"," October 25 2020","Super Kings won by 8 wickets (with 8 balls remaining)","Chennai Super Kings","Royal Challengers Bangalore","25","13","2","2","192.30")
    ,@(" Sharjah"," October 23 2020","Mumbai won by 10 wickets (with 46 balls remaining)","Chennai Super Kings","Mumbai Indians","1","7","0","0","14.28")
    ,@(" Abu Dhabi"," November 01 2020","Super Kings won by 9 wickets (with 7 balls remaining)","Chennai Super Kings","Kings XI Punjab","48","34","4","2","141.17")
    ,@(" Dubai (DSC)"," October 02 2020","Sunrisers won by 7 runs","Chennai Super Kings","Sunrisers Hyderabad","22","19","4","0","115.78")
    ,@(" Sharjah"," September 22 2020","Royals won by 16 runs","Chennai Super Kings","Rajasthan Royals","72","37","1","7","194.59")
    ,@(" Sharjah"," October 17 2020","Capitals won by 5 wickets (with 1 ball remaining)","Chennai Super Kings","Delhi Capitals","58","47","6","2","123.40")
    ,@(" Abu Dhabi"," October 19 2020","Royals won by 7 wickets (with 15 balls remaining)","Chennai Super Kings","Rajasthan Royals","10","9","1","0","111.11")
    ,@(" Dubai (DSC)"," September 25 2020","Capitals won by 44 runs","Chennai Super Kings","Delhi Capitals","43","35","4","0","122.85")
    ,@(" Dubai (DSC)"," October 10 2020","RCB won by 37 runs","Chennai Super Kings","Royal Challengers Bangalore","8","10","0","0","80.00")
    ,@(" Abu Dhabi"," October 07 2020","KKR won by 10 runs","Chennai Super Kings","Kolkata Knight Riders","17","10","3","0","170.00")
)

$startRow = 15
$numericCols = 7,8,9,10,11   # G,H,I,J,K keep "number stored as text"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = "Faf du Plessis" + $nbsp
    for ($c = 0; $c -lt $numericCols.Length; $c++) {
        $col = $numericCols[$c]
        $cell = $ws.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[5 + $c]
    }
}
